$d = $word.ActiveDocument

# Simple whole-document Find/Replace for text that lives in paragraphs
# which contain no "sibling" empty <w:r/> runs (so Word's normal replace
# behavior doesn't disturb any other run in the paragraph).
function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Replace text that sits in a paragraph together with a separate, empty
# <w:r/> run. A plain Range.Text assignment (or Find/Replace) causes this
# engine to rebuild the paragraph's runs and silently merge/drop that
# empty run, which would make the output diverge from the expected
# OOXML. To avoid that, we fetch the paragraph's own WordOpenXML
# (which preserves the exact run layout), patch only the text inside
# it, strip the per-revision/session attributes Word stamps onto the
# <w:p> element when round-tripping, and re-insert that XML in place of
# the matched range - this leaves every other run (including empty
# ones) completely untouched.
function Replace-TextKeepRuns($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    if (-not $rng.Find.Found) {
        return
    }
    $para = $rng.Paragraphs.Item(1)
    $xml = $para.Range.WordOpenXML
    $xml = $xml.Replace($find, $replace)
    $xml = $xml -replace ' w14:paraId="[^"]*"', ''
    $xml = $xml -replace ' w14:textId="[^"]*"', ''
    $xml = $xml -replace ' w:rsidR="[^"]*"', ''
    $xml = $xml -replace ' w:rsidRDefault="[^"]*"', ''
    $xml = $xml -replace ' w:rsidP="[^"]*"', ''
    $xml = $xml -replace ' w:rsidRPr="[^"]*"', ''
    $rng.Delete()
    $rng.InsertXML($xml)
}

# Title (Heading1 at the top, and again as a bold run near the bottom)
Replace-Text "Play Cyberslot Megaclusters for Free - Review" "Play Cyberslot Megaclusters Free - Exciting Features & High RTP"

# "What we like" bullets (each paragraph also holds a leading empty run)
Replace-TextKeepRuns "Innovative gameplay mechanics with the potential for complex grids" "Unique gameplay mechanics with expanding grid"
Replace-TextKeepRuns "Multiplier features can lead to big wins" "Exciting multiplier feature"
Replace-TextKeepRuns "High RTP percentage of 96.36%" "High RTP percentage"
Replace-TextKeepRuns "Futuristic design and theme immerses players in a unique world" "Futuristic theme and immersive soundtrack"

# "What we don't like" bullets (same leading empty run situation)
Replace-TextKeepRuns "Limited symbol design may not be appealing to all players" "Limited symbol variety"
Replace-TextKeepRuns "High volatility can make payouts unpredictable" "Lack of notable sound effects"

# Meta description (italic run)
Replace-Text "Read our review of Cyberslot Megaclusters, an innovative online slot game with unique gameplay mechanics and multiplier features. Play for free here." "Read our review of Cyberslot Megaclusters and discover its unique gameplay and high RTP. Play for free and enjoy the futuristic theme."
